# "fixed horizontal centering on registers"
#
# The workbook's print setup was missing the "Center on page: Horizontally"
# option. Turn it on for the active sheet, which serializes to
# <printOptions horizontalCentered="1"/> in the sheet's OOXML (right before
# <pageMargins .../>), matching the target edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.PageSetup.CenterHorizontally = $true
